$wb = $excel.ActiveWorkbook

# --- ALC row 132 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3699.2646
$ws.Range("I132").Value = 3820.56
$ws.Range("J132").Value = 3362.3333
$ws.Range("K132").Value = 11461.68
$ws.Range("L132").Value = 10086.9999
$ws.Range("M132").Value = -8931.68
$ws.Range("N132").Value = -15146.9999

# --- ALC row 134 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 97986.664
$ws.Range("J134").Value = 97986.664
$ws.Range("L134").Value = 97986.664
$ws.Range("N134").Value = -108126.664

# --- ALC row 139 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 36989.445
$ws.Range("J139").Value = 45129.668
$ws.Range("L139").Value = 45129.668
$ws.Range("N139").Value = -55409.668

# --- ALC row 140 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 48725
$ws.Range("J140").Value = 48725
$ws.Range("L140").Value = 48725
$ws.Range("N140").Value = -59085

# --- ARM row 2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1159.909
$ws.Range("I2").Value = 1108.25
$ws.Range("K2").Value = 1108.25
$ws.Range("M2").Value = -995.25

# --- ARM row 61 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38465156
$ws.Range("I61").Value = 41670336
$ws.Range("K61").Value = 41670336
$ws.Range("M61").Value = -41670124

# --- ARM row 74 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 41672364
$ws.Range("I74").Value = 71432050
$ws.Range("J74").Value = 8805.6
$ws.Range("K74").Value = 71432050
$ws.Range("L74").Value = 8805.6
$ws.Range("M74").Value = -71431176
$ws.Range("N74").Value = -10553.6

# --- ARM row 77 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 41672364
$ws.Range("I77").Value = 71432050
$ws.Range("J77").Value = 8805.6
$ws.Range("K77").Value = 357160250
$ws.Range("L77").Value = 44028
$ws.Range("M77").Value = -357155882
$ws.Range("N77").Value = -52764

# --- ARM row 116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1159.909
$ws.Range("I116").Value = 1108.25
$ws.Range("K116").Value = 1108.25
$ws.Range("M116").Value = 1185.75

# --- ARM row 136 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 38465156
$ws.Range("I136").Value = 41670336
$ws.Range("K136").Value = 125011008
$ws.Range("M136").Value = -125008458

# --- ARM row 139 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 66385.60000000001
$ws.Range("J139").Value = 66385.60000000001
$ws.Range("L139").Value = 66385.60000000001
$ws.Range("N139").Value = -76665.60000000001

# --- BSM row 3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1159.909
$ws.Range("I3").Value = 1108.25
$ws.Range("K3").Value = 1108.25
$ws.Range("M3").Value = -994.25

# --- BSM row 19 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346

# --- BSM row 22 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 197.66667
$ws.Range("I22").Value = 165.57143
$ws.Range("J22").Value = 310
$ws.Range("K22").Value = 165.57143
$ws.Range("L22").Value = 310
$ws.Range("M22").Value = 7.428570000000008
$ws.Range("N22").Value = -656

# --- BSM row 92 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 11850.167
$ws.Range("J92").Value = 11850.167
$ws.Range("L92").Value = 11850.167
$ws.Range("N92").Value = -16842.167

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7250083.5
$ws.Range("I31").Value = 3788.7112
$ws.Range("J31").Value = 333333340
$ws.Range("K31").Value = 3788.7112
$ws.Range("L31").Value = 333333340
$ws.Range("M31").Value = -3493.7112
$ws.Range("N31").Value = -333333930

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7250083.5
$ws.Range("I34").Value = 3788.7112
$ws.Range("J34").Value = 333333340
$ws.Range("K34").Value = 3788.7112
$ws.Range("L34").Value = 333333340
$ws.Range("M34").Value = -3586.7112
$ws.Range("N34").Value = -333333744

# --- CRP row 107 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1041.7894
$ws.Range("I107").Value = 777.2857
$ws.Range("J107").Value = 1782.4
$ws.Range("K107").Value = 777.2857
$ws.Range("L107").Value = 1782.4
$ws.Range("M107").Value = 1142.7143
$ws.Range("N107").Value = -5622.4

# --- CRP row 140 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 42191.355
$ws.Range("J140").Value = 42191.355
$ws.Range("L140").Value = 42191.355
$ws.Range("N140").Value = -52551.355

# --- CUL row 36 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 600
$ws.Range("I36").Value = 200
$ws.Range("J36").Value = 1000
$ws.Range("K36").Value = 600
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = -431
$ws.Range("N36").Value = -3338

# --- CUL row 93 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 4920
$ws.Range("J93").Value = 4920
$ws.Range("L93").Value = 14760
$ws.Range("N93").Value = -18504

# --- GSM row 23 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -1777
$ws.Range("N23").Value = ""

# --- GSM row 80 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13892627
$ws.Range("I80").Value = 23812238
$ws.Range("J80").Value = 5170.6
$ws.Range("K80").Value = 23812238
$ws.Range("L80").Value = 5170.6
$ws.Range("M80").Value = -23811240
$ws.Range("N80").Value = -7166.6

# --- GSM row 83 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 13892627
$ws.Range("I83").Value = 23812238
$ws.Range("J83").Value = 5170.6
$ws.Range("K83").Value = 119061190
$ws.Range("L83").Value = 25853
$ws.Range("M83").Value = -119056198
$ws.Range("N83").Value = -35837

# --- LTW row 22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 965
$ws.Range("I22").Value = 633.3333
$ws.Range("J22").Value = 1296.6666
$ws.Range("K22").Value = 633.3333
$ws.Range("L22").Value = 1296.6666
$ws.Range("M22").Value = -338.3333
$ws.Range("N22").Value = -1886.6666

# --- LTW row 27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 965
$ws.Range("I27").Value = 633.3333
$ws.Range("J27").Value = 1296.6666
$ws.Range("K27").Value = 633.3333
$ws.Range("L27").Value = 1296.6666
$ws.Range("M27").Value = -526.3333
$ws.Range("N27").Value = -1510.6666

# --- LTW row 40 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6337.8335
$ws.Range("I40").Value = 6021.615
$ws.Range("J40").Value = 7160
$ws.Range("K40").Value = 6021.615
$ws.Range("L40").Value = 7160
$ws.Range("M40").Value = -5885.615
$ws.Range("N40").Value = -7432

# --- LTW row 94 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 50958
$ws.Range("J94").Value = 50958
$ws.Range("L94").Value = 50958
$ws.Range("N94").Value = -52310

# --- LTW row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8159.0513
$ws.Range("I132").Value = 4892.8945
$ws.Range("K132").Value = 14678.6835
$ws.Range("M132").Value = -12148.6835

# --- LTW row 136 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 14711697
$ws.Range("I136").Value = 20002308
$ws.Range("J136").Value = 15556.111
$ws.Range("K136").Value = 60006924
$ws.Range("L136").Value = 46668.333
$ws.Range("M136").Value = -60004374
$ws.Range("N136").Value = -51768.333

# --- LTW row 139 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 59918
$ws.Range("J139").Value = 59918
$ws.Range("L139").Value = 59918
$ws.Range("N139").Value = -70198

# --- WVR row 62 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13633.667
$ws.Range("I62").Value = 5333.3335
$ws.Range("J62").Value = 17783.834
$ws.Range("K62").Value = 5333.3335
$ws.Range("L62").Value = 17783.834
$ws.Range("M62").Value = -4709.3335
$ws.Range("N62").Value = -19031.834

# --- WVR row 65 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 13633.667
$ws.Range("I65").Value = 5333.3335
$ws.Range("J65").Value = 17783.834
$ws.Range("K65").Value = 26666.6675
$ws.Range("L65").Value = 88919.17
$ws.Range("M65").Value = -23546.6675
$ws.Range("N65").Value = -95159.17

# --- WVR row 81 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 800
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 800
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 1600
$ws.Range("M81").Value = ""
$ws.Range("N81").Value = -3722

# --- WVR row 84 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 800
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 800
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 8000
$ws.Range("M84").Value = ""
$ws.Range("N84").Value = -18608

# --- WVR row 103 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 27800
$ws.Range("J103").Value = 27800
$ws.Range("L103").Value = 27800
$ws.Range("N103").Value = -30144

# --- WVR row 132 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2236.9443
$ws.Range("I132").Value = 1331.6923
$ws.Range("J132").Value = 4590.6
$ws.Range("K132").Value = 3995.0769
$ws.Range("L132").Value = 13771.8
$ws.Range("M132").Value = -1465.0769
$ws.Range("N132").Value = -18831.8

# --- WVR row 136 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4066.6667
$ws.Range("I136").Value = 3080
$ws.Range("K136").Value = 9240
$ws.Range("M136").Value = -6690
